$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 474.47916
$ws.Range("J17").Value = 474.47916
$ws.Range("L17").Value = 1423.43748
$ws.Range("N17").Value = -1759.43748
$ws.Range("H40").Value = 2614.2856
$ws.Range("I40").Value = 2290
$ws.Range("J40").Value = 3425
$ws.Range("K40").Value = 2290
$ws.Range("L40").Value = 3425
$ws.Range("M40").Value = -2115
$ws.Range("N40").Value = -3775
$ws.Range("H76").Value = 2783.4666
$ws.Range("I76").Value = 2496.121
$ws.Range("J76").Value = 3573.6667
$ws.Range("K76").Value = 2496.121
$ws.Range("L76").Value = 3573.6667
$ws.Range("M76").Value = -2181.121
$ws.Range("N76").Value = -4203.6667
$ws.Range("H79").Value = 2783.4666
$ws.Range("I79").Value = 2496.121
$ws.Range("J79").Value = 3573.6667
$ws.Range("K79").Value = 2496.121
$ws.Range("L79").Value = 3573.6667
$ws.Range("M79").Value = -1404.121
$ws.Range("N79").Value = -5757.6667
$ws.Range("H129").Value = 897.28
$ws.Range("J129").Value = 916.09576
$ws.Range("L129").Value = 2748.28728
$ws.Range("N129").Value = -12748.28728
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1895.8334
$ws.Range("I45").Value = 849
$ws.Range("J45").Value = 2942.6667
$ws.Range("K45").Value = 849
$ws.Range("L45").Value = 2942.6667
$ws.Range("M45").Value = -472
$ws.Range("N45").Value = -3696.6667
$ws.Range("H61").Value = 1488.159
$ws.Range("I61").Value = 1031.8788
$ws.Range("K61").Value = 1031.8788
$ws.Range("M61").Value = -819.8788
$ws.Range("H74").Value = 2506.61
$ws.Range("I74").Value = 2820
$ws.Range("J74").Value = 1139.091
$ws.Range("K74").Value = 2820
$ws.Range("L74").Value = 1139.091
$ws.Range("M74").Value = -1946
$ws.Range("N74").Value = -2887.091
$ws.Range("H77").Value = 2506.61
$ws.Range("I77").Value = 2820
$ws.Range("J77").Value = 1139.091
$ws.Range("K77").Value = 14100
$ws.Range("L77").Value = 5695.455
$ws.Range("M77").Value = -9732
$ws.Range("N77").Value = -14431.455
$ws.Range("H136").Value = 1488.159
$ws.Range("I136").Value = 1031.8788
$ws.Range("K136").Value = 3095.6364
$ws.Range("M136").Value = -545.6363999999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 439.72726
$ws.Range("I94").Value = 388.14285
$ws.Range("J94").Value = 530
$ws.Range("K94").Value = 388.14285
$ws.Range("L94").Value = 530
$ws.Range("M94").Value = 62.85714999999999
$ws.Range("N94").Value = -1432

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1577.0204
$ws.Range("I58").Value = 1069.4
$ws.Range("J58").Value = 2378.5264
$ws.Range("K58").Value = 1069.4
$ws.Range("L58").Value = 2378.5264
$ws.Range("M58").Value = -866.4000000000001
$ws.Range("N58").Value = -2784.5264
$ws.Range("H75").Value = 16000
$ws.Range("I75").Value = 16000
$ws.Range("K75").Value = 16000
$ws.Range("M75").Value = -15002
$ws.Range("H78").Value = 16000
$ws.Range("I78").Value = 16000
$ws.Range("K78").Value = 48000
$ws.Range("M78").Value = -43008
$ws.Range("H105").Value = 698.8889
$ws.Range("I105").Value = 674.75
$ws.Range("J105").Value = 892
$ws.Range("K105").Value = 674.75
$ws.Range("L105").Value = 892
$ws.Range("M105").Value = 1072.25
$ws.Range("N105").Value = -4386
$ws.Range("H132").Value = 1846.3636
$ws.Range("I132").Value = 1658.8182
$ws.Range("J132").Value = 2221.4546
$ws.Range("K132").Value = 4976.4546
$ws.Range("L132").Value = 6664.3638
$ws.Range("M132").Value = -2446.4546
$ws.Range("N132").Value = -11724.3638
$ws.Range("H134").Value = 2167.125
$ws.Range("I134").Value = 1360.6
$ws.Range("J134").Value = 4586.7
$ws.Range("K134").Value = 4081.8
$ws.Range("L134").Value = 13760.1
$ws.Range("M134").Value = -1546.8
$ws.Range("N134").Value = -18830.1
$ws.Range("H136").Value = 1577.0204
$ws.Range("I136").Value = 1069.4
$ws.Range("J136").Value = 2378.5264
$ws.Range("K136").Value = 3208.2
$ws.Range("L136").Value = 7135.5792
$ws.Range("M136").Value = -658.2000000000003
$ws.Range("N136").Value = -12235.5792

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1815345
$ws.Range("I113").Value = 11494720
$ws.Range("J113").Value = 462.125
$ws.Range("K113").Value = 34484160
$ws.Range("L113").Value = 1386.375
$ws.Range("M113").Value = -34481990
$ws.Range("N113").Value = -5726.375

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5267.9546
$ws.Range("I70").Value = 5352.3687
$ws.Range("K70").Value = 5352.3687
$ws.Range("M70").Value = -5082.3687
$ws.Range("H73").Value = 5267.9546
$ws.Range("I73").Value = 5352.3687
$ws.Range("K73").Value = 5352.3687
$ws.Range("M73").Value = -4416.3687
$ws.Range("H97").Value = 629.7273
$ws.Range("I97").Value = 641.7
$ws.Range("J97").Value = 510
$ws.Range("K97").Value = 641.7
$ws.Range("L97").Value = 510
$ws.Range("M97").Value = -145.7
$ws.Range("N97").Value = -1502

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 734.7692
$ws.Range("I22").Value = 592
$ws.Range("J22").Value = 1210.6666
$ws.Range("K22").Value = 592
$ws.Range("L22").Value = 1210.6666
$ws.Range("M22").Value = -297
$ws.Range("N22").Value = -1800.6666
$ws.Range("H27").Value = 734.7692
$ws.Range("I27").Value = 592
$ws.Range("J27").Value = 1210.6666
$ws.Range("K27").Value = 592
$ws.Range("L27").Value = 1210.6666
$ws.Range("M27").Value = -485
$ws.Range("N27").Value = -1424.6666
$ws.Range("H46").Value = 1220
$ws.Range("I46").Value = 960
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 960
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -772
$ws.Range("N46").Value = -2376
$ws.Range("H55").Value = 810.26086
$ws.Range("I55").Value = 804.3158
$ws.Range("J55").Value = 838.5
$ws.Range("K55").Value = 804.3158
$ws.Range("L55").Value = 838.5
$ws.Range("M55").Value = -631.3158
$ws.Range("N55").Value = -1184.5
$ws.Range("H93").Value = 11405
$ws.Range("I93").Value = 26600.5
$ws.Range("J93").Value = 1274.6666
$ws.Range("K93").Value = 26600.5
$ws.Range("L93").Value = 1274.6666
$ws.Range("M93").Value = -25352.5
$ws.Range("N93").Value = -3770.6666
$ws.Range("H132").Value = 2278.27
$ws.Range("I132").Value = 1897.8375
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 5693.512500000001
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -3163.512500000001
$ws.Range("N132").Value = -16460
$ws.Range("H133").Value = 45467.375
$ws.Range("J133").Value = 45467.375
$ws.Range("L133").Value = 45467.375
$ws.Range("N133").Value = -50527.375
$ws.Range("H136").Value = 2017.6615
$ws.Range("I136").Value = 1624.5217
$ws.Range("J136").Value = 2969.4736
$ws.Range("K136").Value = 4873.5651
$ws.Range("L136").Value = 8908.4208
$ws.Range("M136").Value = -2323.5651
$ws.Range("N136").Value = -14008.4208

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 33271
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 33271
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 33271
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -34893
$ws.Range("H71").Value = 33271
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 33271
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 99813
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -107925
$ws.Range("H96").Value = 1688.6923
$ws.Range("I96").Value = 1060.75
$ws.Range("J96").Value = 1967.7778
$ws.Range("K96").Value = 1060.75
$ws.Range("L96").Value = 1967.7778
$ws.Range("M96").Value = 312.25
$ws.Range("N96").Value = -4713.7778
$ws.Range("H107").Value = 324.5
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 312.9091
$ws.Range("K107").Value = 1050
$ws.Range("L107").Value = 938.7273
$ws.Range("M107").Value = 870
$ws.Range("N107").Value = -4778.7273
$ws.Range("H132").Value = 2341.204
$ws.Range("I132").Value = 2364.1614
$ws.Range("J132").Value = 2301.6667
$ws.Range("K132").Value = 7092.4842
$ws.Range("L132").Value = 6905.000100000001
$ws.Range("M132").Value = -4562.4842
$ws.Range("N132").Value = -11965.0001
$ws.Range("H138").Value = 47592
$ws.Range("J138").Value = 47592
$ws.Range("L138").Value = 47592
$ws.Range("N138").Value = -57872
